$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
# Target widths: col A -> 22, col B -> 12 (ColumnWidth setter adds ~0.8333 padding)
$ws.Range("A1").EntireColumn.ColumnWidth = 21.1666666666667
$ws.Range("B1").EntireColumn.ColumnWidth = 11.1666666666667

# --- Sheet view ---
$ws.Application.ActiveWindow.ScrollRow = 64
$ws.Range("D75").Select()

# --- New content: string cells first, in authoring order, so shared-string
#     table indices line up with the target workbook ---
$ws.Range("A72").Value = "Robot in Middle of Unit"
$ws.Range("A73").Value = "Left Sensor"
$ws.Range("A81").Value = "Right Sensor"
$ws.Range("C73").Value = "STDEV"
$ws.Range("D73").Value = "CM"
$ws.Range("B73").Value = "AVG"
$ws.Range("B81").Value = "AVG"
$ws.Range("C81").Value = "STDEV"

# --- Alignment / number-format styles, in the order they first appear so
#     new cellXfs entries land at the same indices as the target file ---
# Left-aligned (style index 8 in target)
$ws.Range("A81").HorizontalAlignment = -4131
# Right-aligned (style index 9 in target)
$ws.Range("A82:A87").HorizontalAlignment = -4152
# Integer number format + centered (style index 10 in target)
$rng10 = $ws.Range("B74:C82")
$rng10.NumberFormat = "0"
$rng10.HorizontalAlignment = -4108

# --- Numeric data: left sensor trial block ---
$ws.Range("A74").Value = 534
$ws.Range("A75").Value = 531
$ws.Range("A76").Value = 540
$ws.Range("A77").Value = 551
$ws.Range("A78").Value = 535
$ws.Range("A79").Value = 565
$ws.Range("D74").Value = 7.348

# --- Numeric data: right sensor trial block ---
$ws.Range("A82").Value = 409
$ws.Range("A83").Value = 447
$ws.Range("A84").Value = 432
$ws.Range("A85").Value = 418
$ws.Range("A86").Value = 445
$ws.Range("A87").Value = 408

# --- Formulas ---
$ws.Range("B74").Formula = "=AVERAGE(A74:A79)"
$ws.Range("C74").Formula = "=STDEV.P(A74:A79)"
$ws.Range("B82").Formula = "=AVERAGE(A82:A87)"
$ws.Range("C82").Formula = "=STDEV.P(A82:A87)"
